# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) throughout the resume body.
#
# Strategy: for each paragraph that needs highlighting, first Find.Execute
# the paragraph's full (pre-edit) text against $d.Content - this narrows
# the Range to exactly that paragraph (and is immune to the many duplicate
# numbers/percentages that occur elsewhere in the document). Then, within
# that narrowed Range, Find.Execute each individual metric substring with
# Replacement formatting (Bold + dark slate color 2C3E50) and the *same*
# text as the replacement - this is a pure reformat, so Word splits the
# run around the match exactly like the diff shows, without altering the
# surrounding wording.

$d = $word.ActiveDocument

# wdColor integer for hex 2C3E50 (VBA RGB() packs as R + G*256 + B*65536)
$metricColor = 5258796

function Set-MetricBold($ParagraphText, $Metrics) {
    $whole = $d.Content
    $found = $whole.Find.Execute($ParagraphText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $ParagraphText"
        return
    }

    $paraStart = $whole.Start
    $paraEnd = $whole.End

    foreach ($metric in $Metrics) {
        $scoped = $d.Range($paraStart, $paraEnd)
        $f = $scoped.Find
        $f.ClearFormatting()
        $f.Replacement.ClearFormatting()
        $f.Replacement.Font.Bold = $true
        $f.Replacement.Font.Color = $metricColor
        $f.Forward = $true
        $f.Wrap = 0
        $ok = $f.Execute($metric, $false, $false, $false, $false, $false, $true, 1, $false, $metric, 2)
        if (-not $ok) {
            Write-Output "  metric not matched in scope: $metric (paragraph: $ParagraphText)"
        }
    }
}

# --- Siege Analytics bullets ---

Set-MetricBold `
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" `
    @("23%", "64%")

Set-MetricBold `
    "• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes" `
    @("±4.2%", "±2.1%", "71%", "87%")

Set-MetricBold `
    "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" `
    @("73.5%", "`$4.7M")

Set-MetricBold `
    "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" `
    @("`$2")

# --- Helm/Murmuration bullet ---

Set-MetricBold `
    "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%" `
    @("57%")

# --- KEY ACHIEVEMENTS AND IMPACT bullets ---

Set-MetricBold `
    "• 178% accuracy improvement in racial classification algorithms" `
    @("178%")

Set-MetricBold `
    "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%" `
    @("73.5%")

Set-MetricBold `
    "• `$4.7M savings enabled nonprofit access" `
    @("`$4.7M")

Set-MetricBold `
    "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations" `
    @("12,847")

Set-MetricBold `
    "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%" `
    @("±4.2%", "±2.1%")

Set-MetricBold `
    "• Increased voter turnout prediction accuracy from 71% to 87%" `
    @("71%", "87%")

Write-Output "Done."
